# Update column F (dSF) values for several rows to reflect repulled data /
# recalculated means, per commit message "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    6  = -4
    9  = 8
    10 = -10
    11 = -9
    12 = 9
    13 = -3
    14 = 2
    16 = -2
    18 = -6
    19 = 9
    22 = -3
    23 = -8
    24 = -1
    25 = -2
    26 = -4
    28 = -1
    30 = 4
    32 = -5
    37 = 1
    38 = -7
    39 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
